# Broker Files can now be modified
#
# The workbook previously shipped with workbook-level protection enabled
# (an empty <workbookProtection/> record) which is being lifted so the
# broker chart can be edited going forward. Two "Settled" installments are
# corrected to "Paid", and the active selection/view is moved to the last
# edited cell (D24) to mirror where the editor was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lift workbook protection (structure/windows) so the sheet can be edited.
$wb.Unprotect()

# Correct the payment status for the last two installments.
$ws.Range("C25").Value = "Paid"
$ws.Range("C26").Value = "Paid"

# Leave the selection where the edit happened.
$ws.Range("D24").Select()
